$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update odds columns ---
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.4
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("V2").Value = 1.67

# --- Row 3: update odds columns ---
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("V3").Value = 1.62

# --- Row 4 (INDIA - ISL / Hyderabad vs Odisha FC) is removed entirely.
# Deleting the whole row shifts rows 5 and 6 up to become rows 4 and 5,
# and the sheet's used range shrinks accordingly.
$ws.Rows.Item(4).Delete()

# --- The new row 4 (previously row 5) has two odds tweaked further ---
$ws.Range("K4").Value = 1.95
$ws.Range("R4").Value = 1.53
